$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.793.06"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.613.56"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'539.48"
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("D6").Value = "'142.50"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'6.48"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'0.336"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "3.065.16"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "59.701.18"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "'20.75"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'0.0000134"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "2.605.04"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "'342.99"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'4.38"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "'10.17"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'67.78"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "'0.411"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'7.26"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'1.68"
$ws.Range("E30").Value = "  +6.28%  "
$ws.Range("D31").Value = "'5.88"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'18.93"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'149.82"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'4.00"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'1.12"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'0.854"
$ws.Range("E36").Value = "  +5.40%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("D38").Value = "'0.833"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "'3.56"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'274.85"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'0.598"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.0959"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'10.73"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'0.0525"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").Value = "1.961.43"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'18.57"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0224"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").Value = "'4.54"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "'112.61"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("E51").Value = "  +0.54%  "
